$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Execution metrics section updates
# B5: Compilation success -> "no"
$ws.Range("B5").Value = "no"
# C5: Note -> "Called wrong method"
$ws.Range("C5").Value = "Called wrong method"

# B6: Runtime without error -> cleared (was "yes")
$ws.Range("B6").Value = ""

# B7: Assertion validity -> cleared (was "yes")
$ws.Range("B7").Value = ""
# C7: Note -> cleared (was "Test fail, same as baseline")
$ws.Range("C7").Value = ""

# Code BLEU updated score
$ws.Range("B12").Value = 0.2690207408326153
$ws.Range("C12").Value = "{'codebleu': 0.26902074083261535, 'ngram_match_score': 0.08993399093494622, 'weighted_ngram_match_score': 0.10530432684498735, 'syntax_match_score': 0.6102564102564103, 'dataflow_match_score': 0.27058823529411763}"

# Update the active selection on the sheet to B6 (matches the recorded cursor position)
$ws.Range("B6").Select()

$wb.Save()
